# Auto-generated edit script: apply scheduled-runner price/profit updates
# to the Alexander_Profits workbook (8 Leve-crafting sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4879.3184
$ws.Range("J74").Value = 3045.4546
$ws.Range("L74").Value = 3045.4546
$ws.Range("N74").Value = -4917.4546
$ws.Range("H77").Value = 4879.3184
$ws.Range("J77").Value = 3045.4546
$ws.Range("L77").Value = 15227.273
$ws.Range("N77").Value = -24587.273
$ws.Range("H116").Value = 4042.561
$ws.Range("I116").Value = 3961.818
$ws.Range("J116").Value = 4136.0527
$ws.Range("K116").Value = 3961.818
$ws.Range("L116").Value = 4136.0527
$ws.Range("M116").Value = -519.8180000000002
$ws.Range("N116").Value = -11020.0527
$ws.Range("H132").Value = 2020.3971
$ws.Range("I132").Value = 1322
$ws.Range("J132").Value = 4975.154
$ws.Range("K132").Value = 3966
$ws.Range("L132").Value = 14925.462
$ws.Range("M132").Value = -1436
$ws.Range("N132").Value = -19985.462
$ws.Range("H135").Value = 28772.277
$ws.Range("I135").Value = 34158.133
$ws.Range("J135").Value = 1843
$ws.Range("K135").Value = 307423.197
$ws.Range("L135").Value = 16587
$ws.Range("M135").Value = -304888.197
$ws.Range("N135").Value = -21657
$ws.Range("H137").Value = 3126134
$ws.Range("I137").Value = 1429686.5
$ws.Range("K137").Value = 4289059.5
$ws.Range("M137").Value = -4286509.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2510.97
$ws.Range("I32").Value = 2256.9238
$ws.Range("J32").Value = 5432.5
$ws.Range("K32").Value = 2256.9238
$ws.Range("L32").Value = 5432.5
$ws.Range("M32").Value = -1969.9238
$ws.Range("N32").Value = -6006.5
$ws.Range("H74").Value = 855.5111000000001
$ws.Range("I74").Value = 877.2381
$ws.Range("J74").Value = 551.3333
$ws.Range("K74").Value = 877.2381
$ws.Range("L74").Value = 551.3333
$ws.Range("M74").Value = -3.238100000000031
$ws.Range("N74").Value = -2299.3333
$ws.Range("H77").Value = 855.5111000000001
$ws.Range("I77").Value = 877.2381
$ws.Range("J77").Value = 551.3333
$ws.Range("K77").Value = 4386.190500000001
$ws.Range("L77").Value = 2756.6665
$ws.Range("M77").Value = -18.19050000000061
$ws.Range("N77").Value = -11492.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2109.0417
$ws.Range("I86").Value = 2271.8235
$ws.Range("J86").Value = 1713.7142
$ws.Range("K86").Value = 2271.8235
$ws.Range("L86").Value = 1713.7142
$ws.Range("M86").Value = -1148.8235
$ws.Range("N86").Value = -3959.7142
$ws.Range("H89").Value = 2109.0417
$ws.Range("I89").Value = 2271.8235
$ws.Range("J89").Value = 1713.7142
$ws.Range("K89").Value = 11359.1175
$ws.Range("L89").Value = 8568.571
$ws.Range("M89").Value = -5743.1175
$ws.Range("N89").Value = -19800.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1655.093
$ws.Range("I31").Value = 1382.1428
$ws.Range("J31").Value = 2849.25
$ws.Range("K31").Value = 1382.1428
$ws.Range("L31").Value = 2849.25
$ws.Range("M31").Value = -1087.1428
$ws.Range("N31").Value = -3439.25
$ws.Range("H34").Value = 1655.093
$ws.Range("I34").Value = 1382.1428
$ws.Range("J34").Value = 2849.25
$ws.Range("K34").Value = 1382.1428
$ws.Range("L34").Value = 2849.25
$ws.Range("M34").Value = -1180.1428
$ws.Range("N34").Value = -3253.25
$ws.Range("H58").Value = 768.55
$ws.Range("I58").Value = 796.6445
$ws.Range("J58").Value = 684.26666
$ws.Range("K58").Value = 796.6445
$ws.Range("L58").Value = 684.26666
$ws.Range("M58").Value = -593.6445
$ws.Range("N58").Value = -1090.26666
$ws.Range("H62").Value = 3158.1177
$ws.Range("I62").Value = 3045.8667
$ws.Range("K62").Value = 3045.8667
$ws.Range("M62").Value = -2421.8667
$ws.Range("H65").Value = 3158.1177
$ws.Range("I65").Value = 3045.8667
$ws.Range("K65").Value = 15229.3335
$ws.Range("M65").Value = -12109.3335
$ws.Range("H134").Value = 5352.8184
$ws.Range("I134").Value = 5747.107
$ws.Range("J134").Value = 3144.8
$ws.Range("K134").Value = 17241.321
$ws.Range("L134").Value = 9434.400000000001
$ws.Range("M134").Value = -14706.321
$ws.Range("N134").Value = -14504.4
$ws.Range("H136").Value = 768.55
$ws.Range("I136").Value = 796.6445
$ws.Range("J136").Value = 684.26666
$ws.Range("K136").Value = 2389.9335
$ws.Range("L136").Value = 2052.79998
$ws.Range("M136").Value = 160.0664999999999
$ws.Range("N136").Value = -7152.79998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 916.1194
$ws.Range("I131").Value = 531
$ws.Range("J131").Value = 975.87933
$ws.Range("K131").Value = 1593
$ws.Range("L131").Value = 2927.63799
$ws.Range("M131").Value = 3447
$ws.Range("N131").Value = -13007.63799
$ws.Range("H136").Value = 4068.349
$ws.Range("I136").Value = 2342
$ws.Range("J136").Value = 4394.0757
$ws.Range("K136").Value = 7026
$ws.Range("L136").Value = 13182.2271
$ws.Range("M136").Value = -1926
$ws.Range("N136").Value = -23382.2271

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1325471.5
$ws.Range("I80").Value = 2627.9285
$ws.Range("J80").Value = 2750072.2
$ws.Range("K80").Value = 2627.9285
$ws.Range("L80").Value = 2750072.2
$ws.Range("M80").Value = -1629.9285
$ws.Range("N80").Value = -2752068.2
$ws.Range("H83").Value = 1325471.5
$ws.Range("I83").Value = 2627.9285
$ws.Range("J83").Value = 2750072.2
$ws.Range("K83").Value = 13139.6425
$ws.Range("L83").Value = 13750361
$ws.Range("M83").Value = -8147.6425
$ws.Range("N83").Value = -13760345

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3833.3333
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 4250
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 4250
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -4474
$ws.Range("H40").Value = 2117.818
$ws.Range("I40").Value = 1780
$ws.Range("J40").Value = 2399.3333
$ws.Range("K40").Value = 1780
$ws.Range("L40").Value = 2399.3333
$ws.Range("M40").Value = -1644
$ws.Range("N40").Value = -2671.3333
$ws.Range("H136").Value = 1342.1177
$ws.Range("I136").Value = 1144.7222
$ws.Range("J136").Value = 2103.5
$ws.Range("K136").Value = 3434.1666
$ws.Range("L136").Value = 6310.5
$ws.Range("M136").Value = -884.1665999999996
$ws.Range("N136").Value = -11410.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2247.5962
$ws.Range("I132").Value = 2210.681
$ws.Range("J132").Value = 2594.6
$ws.Range("K132").Value = 6632.043
$ws.Range("L132").Value = 7783.799999999999
$ws.Range("M132").Value = -4102.043
$ws.Range("N132").Value = -12843.8
